# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Kazajistan (row 62): Muertes hoy 31 -> 35 ---
$ws.Range("F62").Value = 35

# --- Uzbekistan (row 68): Casos totales / Nuevos casos / Recuperados ---
$ws.Range("B68").Value = 1887
$ws.Range("C68").Value = 18
$ws.Range("E68").Value = 1090

# --- Re-rank around Georgia/Malta/Jordania/Consejo Danes/Somalia/Taiwan ---
# "Consejo Danes para los Refugiados" gets updated figures that move it
# above Malta and Jordania (rows 109-111), shifting them down one row.
$ws.Range("A109").Value = "Consejo Danes para los Refugiados"
$ws.Range("B109").Value = 459
$ws.Range("C109").Value = 17
$ws.Range("D109").Value = 50
$ws.Range("E109").Value = 381
$ws.Range("F109").Value = 0
$ws.Range("G109").Value = 0
$ws.Range("H109").Value = 28

$ws.Range("A110").Value = "Malta"
$ws.Range("B110").Value = 448
$ws.Range("C110").Value = 0
$ws.Range("D110").Value = 282
$ws.Range("E110").Value = 162
$ws.Range("F110").Value = 2
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 4

$ws.Range("A111").Value = "Jordania"
$ws.Range("B111").Value = 447
$ws.Range("C111").Value = 0
$ws.Range("D111").Value = 337
$ws.Range("E111").Value = 103
$ws.Range("F111").Value = 5
$ws.Range("G111").Value = 0
$ws.Range("H111").Value = 7

# Row 112 (Somalia) is unchanged.

# --- Taiwan (row 113): Casos activos / Recuperados ---
$ws.Range("D113").Value = 290
$ws.Range("E113").Value = 133

# --- Zimbabue (row 177): Casos activos / Recuperados ---
$ws.Range("D177").Value = 5
$ws.Range("E177").Value = 22
